$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.082.85'
$ws.Range('E2').Value = '  +0.39%  '
$ws.Range('D3').Value = '1.850.43'
$ws.Range('E3').Value = '  +1.66%  '
$ws.Range('E4').Value = '  +0.52%  '
$ws.Range('D5').Formula = "'237.44"
$ws.Range('E5').Value = '  +3.05%  '
$ws.Range('E6').Value = '  +0.71%  '
$ws.Range('E7').Value = '  +0.30%  '
$ws.Range('D8').Formula = "'42.46"
$ws.Range('E8').Value = '  +5.94%  '
$ws.Range('E9').Value = '  +1.27%  '
$ws.Range('D10').Formula = "'0.0691"
$ws.Range('E10').Value = '  +1.19%  '
$ws.Range('E11').Value = '  +0.08%  '
$ws.Range('D12').Value = '2.118.28'
$ws.Range('E12').Value = '  +1.73%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.855.87'
$ws.Range('E13').Value = '  +1.99%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').Formula = "'11.38"
$ws.Range('E14').Value = '  +0.92%  '
$ws.Range('E15').Value = '  +0.84%  '
$ws.Range('E16').Value = '  +3.10%  '
$ws.Range('D17').Value = '35.029.07'
$ws.Range('E17').Value = '  +0.37%  '
$ws.Range('D18').Formula = "'69.98"
$ws.Range('E18').Value = '  +0.38%  '
$ws.Range('D19').Value = '0.0₃0791'
$ws.Range('E19').Value = '  +0.81%  '
$ws.Range('D20').Formula = "'240.46"
$ws.Range('E20').Value = '  -0.20%  '
$ws.Range('D21').Formula = "'12.14"
$ws.Range('E21').Value = '  +0.92%  '
$ws.Range('E22').Value = '  +1.58%  '
$ws.Range('E23').Value = '  +0.28%  '
$ws.Range('E24').Value = '  -0.37%  '
$ws.Range('D25').Formula = "'169.77"
$ws.Range('E25').Value = '  -2.19%  '
$ws.Range('E26').Value = '  +3.02%  '
$ws.Range('D27').Formula = "'1.83"
$ws.Range('E27').Value = '  +21.08%  '
$ws.Range('D28').Formula = "'17.59"
$ws.Range('E28').Value = '  +1.49%  '
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('E30').Value = '  +0.51%  '
$ws.Range('D31').Formula = "'0.0552"
$ws.Range('E31').Value = '  +0.67%  '
$ws.Range('D32').Formula = "'3.98"
$ws.Range('E32').Value = '  -0.50%  '
$ws.Range('D33').Formula = "'4.01"
$ws.Range('E33').Value = '  +1.32%  '
$ws.Range('E34').Value = '  +25.45%  '
$ws.Range('D35').Formula = "'1.99"
$ws.Range('E35').Value = '  +9.41%  '
$ws.Range('D36').Formula = "'0.796"
$ws.Range('E36').Value = '  +15.16%  '
$ws.Range('E37').Value = '  +4.14%  '
$ws.Range('E38').Value = '  +9.31%  '
$ws.Range('E39').Value = '  +4.22%  '
$ws.Range('D40').Formula = "'90.12"
$ws.Range('E40').Value = '  -2.71%  '
$ws.Range('D41').Value = '1.341.98'
$ws.Range('E41').Value = '  +0.19%  '
$ws.Range('B42').Value = 'InjectiveProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D42').Formula = "'14.86"
$ws.Range('E42').Value = '  +2.33%  '
$ws.Range('B43').Value = 'Gas'
$ws.Range('C43').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D43').Formula = "'12.96"
$ws.Range('E43').Value = '  +52.38%  '
$ws.Range('D44').Formula = "'2.30"
$ws.Range('E44').Value = '  +1.12%  '
$ws.Range('D45').Formula = "'2.44"
$ws.Range('E45').Value = '  +1.24%  '
$ws.Range('D46').Formula = "'2.74"
$ws.Range('E46').Value = '  -0.59%  '
$ws.Range('D47').Formula = "'0.0552"
$ws.Range('E47').Value = '  +6.09%  '
$ws.Range('D48').Formula = "'6.46"
$ws.Range('E48').Value = '  +4.11%  '
$ws.Range('D49').Value = '2.038.32'
$ws.Range('E49').Value = '  +2.06%  '
$ws.Range('D50').Formula = "'0.0674"
$ws.Range('E50').Value = '  +1.57%  '
$ws.Range('E51').Value = '  +0.22%  '
